# TeamANSReport.xlsx — "Add files via upload"
#
# The backlog's Sprint-2 section gains two new user stories (US17, US18)
# that were picked up by a new team member "AA", and the two Sprint-2
# stories already "Coding" (US07 / US10) are marked "Done". The same two
# new stories are also added to the Sprint2 worksheet with their
# estimates. Finally Sprint2 becomes the active/selected sheet (it was
# Backlog before).

$wb = $excel.ActiveWorkbook

# --- Backlog sheet: mark the two in-progress Sprint-1 items as Done ---
$backlog = $wb.Worksheets.Item("Backlog")
$backlog.Range("E8").Value = "Done"
$backlog.Range("E9").Value = "Done"

# --- Backlog sheet: append the two new Sprint-2 backlog items ---
$backlog.Range("A14").Value = 2
$backlog.Range("B14").Value = "US18"
$backlog.Range("C14").Value = "Siblings should not marry"
$backlog.Range("D14").Value = "AA"
$backlog.Range("E14").Value = "Coding"

$backlog.Range("A15").Value = 2
$backlog.Range("B15").Value = "US17"
$backlog.Range("C15").Value = "No marriages to descendants"
$backlog.Range("D15").Value = "AA"
$backlog.Range("E15").Value = "Coding"

$null = $backlog.Range("B14:E15").Select()

# --- Sprint2 sheet: add the same two stories with their estimates ---
$sprint2 = $wb.Worksheets.Item("Sprint2")
$sprint2.Range("A5").Value = "US18"
$sprint2.Range("B5").Value = "Siblings should not marry"
$sprint2.Range("C5").Value = "AA"
$sprint2.Range("D5").Value = "Coding"
$sprint2.Range("E5").Value = 80
$sprint2.Range("F5").Value = 100

$sprint2.Range("A6").Value = "US17"
$sprint2.Range("B6").Value = "No marriages to descendants"
$sprint2.Range("C6").Value = "AA"
$sprint2.Range("D6").Value = "Coding"
$sprint2.Range("E6").Value = 80
$sprint2.Range("F6").Value = 100

# Sprint2 becomes the active sheet/tab, with F11 as the resting selection.
$null = $sprint2.Activate()
$null = $sprint2.Range("F11").Select()
